$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 364, shifting existing rows 364-388 down to 365-389
$ws.Range("A364").EntireRow.Insert()

# Populate the newly inserted row 364 with the new weekly price record
$ws.Range("A364").Value = 7
$ws.Range("B364").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C364").Value = "Ñuble"
$ws.Range("D364").Value = 45267
$ws.Range("E364").Value = 16
$ws.Range("F364").Value = 100112032
$ws.Range("G364").Value = "Zapallo italiano"
$ws.Range("H364").Value = "Sin especificar"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 160
$ws.Range("K364").Value = 10000
$ws.Range("L364").Value = 11000
$ws.Range("M364").Value = 10625
$ws.Range("N364").Value = "$/caja 50 unidades"
$ws.Range("O364").Value = "Región del Maule"
$ws.Range("P364").Value = 212
$ws.Range("Q364").Value = 50
$ws.Range("R364").Value = "Hortaliza"
